$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.231.60"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "2.264.36"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.60"

$ws.Range("D6").Value = "128.78"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("D9").Value = "0.0953"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("D11").Value = "0.334"
$ws.Range("E11").Value = "  +2.68%  "

$ws.Range("D12").Value = "4.83"
$ws.Range("E12").Value = "  +3.91%  "

$ws.Range("E13").Value = "  +5.31%  "

$ws.Range("D14").Value = "2.665.02"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").Value = "54.208.08"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").Value = "2.276.39"
$ws.Range("E17").Value = "  -4.15%  "

$ws.Range("D18").Value = "10.22"
$ws.Range("E18").Value = "  +1.91%  "

$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "303.12"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "60.68"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("E24").Value = "  -2.18%  "

$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").Value = "7.26"
$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("D27").Value = "171.72"
$ws.Range("E27").Value = "  +1.91%  "

$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").Value = "5.96"
$ws.Range("E29").Value = "  +1.61%  "

$ws.Range("D30").Value = "0.0₃0689"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.80"
$ws.Range("E33").Value = "  +0.46%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").Value = "0.943"
$ws.Range("E35").Value = "  +4.00%  "

$ws.Range("E36").Value = "  +0.41%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "0.375"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").Value = "4.81"
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").Value = "124.62"
$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("D43").Value = "0.0491"
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").Value = "0.0894"
$ws.Range("E44").Value = "  +0.79%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").Value = "241.01"
$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("D47").Value = "0.373"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.10"
$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("E51").Value = "  -0.42%  "
